# Add team record (Wins/Losses/Ties) columns to the CHC_2017 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add Wins / Losses / Ties, matching existing header style ---
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-49): every player on this roster shares the team's 92-70-0 record ---
for ($r = 2; $r -le 49; $r++) {
    $ws.Cells.Item($r, 30).Value = 92
    $ws.Cells.Item($r, 31).Value = 70
    $ws.Cells.Item($r, 32).Value = 0
}
